$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Update F2: change item code from TB7SX6CC to TB7SX1CC
$ws.Range("F2").Value = "TB7SX1CC"

# Add new row 4: F4 = TB7SX6CC, same style as F2/F3 (text format)
$ws.Range("F4").Value = "TB7SX6CC"
$ws.Range("F4").NumberFormat = $ws.Range("F3").NumberFormat
